# Updating Self Learning Test Cases: replace TC3 (Nmap CLI) content with
# TC2 (Zenmap GUI) content, rename the sheet, and clear out the rows/
# cells that no longer apply.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet name: "Self Learning TC3" -> "Self Learning TC2" ---
$ws.Name = "Self Learning TC2"

# --- A1: drop the old AC4 heading text entirely, restore natural height ---
$ws.Range("A1").ClearContents()
$ws.Rows(1).AutoFit()

# --- Row 11: new TC2 scenario / first step; clear the old expected result ---
$ws.Range("A11").Value = "TC3: Use quick-scan on NMAP to scan the websites Ipaddress and locate any open ports"
$ws.Range("B11").Value = '1.Initialize the "Zenmap" application from your desktop'
$ws.Range("C11").ClearContents()

# --- Row 12: new step 2 + expected result ---
$ws.Range("B12").Value = '2. In the "Target:" box on top type in the web address "167.71.183.120"'
$ws.Range("C12").Value = "Nmap should show the Ipaddress you entered into the target box as well as update the command box underneath showing you the address there as well. "

# --- Row 13: new step 3 + expected result ---
$ws.Range("B13").Value = '3. On the right side of the application ensure "Quick Scan" is selected as the profile option and then click the "Scan" commmand. '
$ws.Range("C13").Value = "Nmap should populate a list of the open ports located through the Ipaddress you scanned. Take note of which ports are open. Port 443/tcp mysql is a commonly targeted port by hackers. You should take measure to ensure that these ports is properly secured. "

# --- Row 14: no longer used (was step 4 / its result) - clear it out and shrink it back down ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows(14).RowHeight = 15.6

# --- Selection moves to A13, scrolled so row 3 is at the top of the view ---
$ws.Range("A13").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
